$p = $ppt.ActivePresentation
$props = @("Package","PackagePart","Parts","OpenXmlPackage","FullName","Path")
foreach ($pr in $props) {
  try {
    $v = $p.$pr
    Write-Host "$pr => $v"
  } catch {
    Write-Host "$pr => ERROR $_"
  }
}
